# Latest Updated Script for Parallel Execution
#
# Reproduces the newest "test run" cycle of the ExamCenterDetails
# automation workbook (STAGE sheet): the data row (row 2) gets a freshly
# generated Exam-Center Location plus matching Exam/Schedule names, and
# the three LastName cells (E3:E5) get freshly generated numeric ids.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("STAGE")

# --- newest batch of generated test data -----------------------------
$newLocation     = "ECLocation21953"
$newExamName     = "FPK12Exam38575"
$newScheduleName = "FPK12Schedule10282"
$newLastName3    = "87737"
$newLastName4    = "87927"
$newLastName5    = "37737"

# --- row 2: Location / ExamName / ScheduleName ------------------------
$ws.Range("A2").Value = $newLocation
$ws.Range("H2").Value = $newExamName
$ws.Range("I2").Value = $newScheduleName

# --- rows 3-5: LastName (kept as text, like the original numeric-looking
#     ids, via a leading apostrophe so Excel doesn't coerce them to numbers)
$ws.Range("E3").Value = "'" + $newLastName3
$ws.Range("E4").Value = "'" + $newLastName4
$ws.Range("E5").Value = "'" + $newLastName5
